# CreateDataConfigTest.xlsx — "updated till mfg and subset"
#
# 1. createChildCategories!F2:F17 — rename the "ATChildCAT#" values that were
#    used as placeholders to the real "AutomationTestCat#" values (4 distinct
#    strings, 4 rows each).
# 2. createChildCategories column F widened to fit the new (longer) text and
#    the sheet becomes the active tab with a new selection.
# 3. createCategories loses the "active tab" flag and gets a fresh selection.
# 4. The workbook-level view scrolls the tab strip and switches the active
#    sheet to createChildCategories (last sheet).

$wb = $excel.ActiveWorkbook

# --- createChildCategories: update the F-column values -------------------
$childCatWs = $wb.Worksheets.Item("createChildCategories")

$childCatWs.Range("F2").Value  = "AutomationTestCat1"
$childCatWs.Range("F3").Value  = "AutomationTestCat1"
$childCatWs.Range("F4").Value  = "AutomationTestCat1"
$childCatWs.Range("F5").Value  = "AutomationTestCat1"
$childCatWs.Range("F6").Value  = "AutomationTestCat2"
$childCatWs.Range("F7").Value  = "AutomationTestCat2"
$childCatWs.Range("F8").Value  = "AutomationTestCat2"
$childCatWs.Range("F9").Value  = "AutomationTestCat2"
$childCatWs.Range("F10").Value = "AutomationTestCat3"
$childCatWs.Range("F11").Value = "AutomationTestCat3"
$childCatWs.Range("F12").Value = "AutomationTestCat3"
$childCatWs.Range("F13").Value = "AutomationTestCat3"
$childCatWs.Range("F14").Value = "AutomationTestCat4"
$childCatWs.Range("F15").Value = "AutomationTestCat4"
$childCatWs.Range("F16").Value = "AutomationTestCat4"
$childCatWs.Range("F17").Value = "AutomationTestCat4"

# Widen column F to fit the longer strings (matches the other bestFit columns).
$childCatWs.Columns.Item(6).ColumnWidth = 18.6

# --- createCategories: drop the "active" flag, move the selection --------
$categoriesWs = $wb.Worksheets.Item("createCategories")
$categoriesWs.Range("J24").Select()

# --- createChildCategories becomes the active sheet/tab ------------------
$childCatWs.Activate()
$childCatWs.Range("G12").Select()

# Scroll the workbook tab strip so the first visible tab is index 2
# (0-based), mirroring the new firstSheet="2" view setting.
$window = $excel.ActiveWindow
$window.ScrollWorkbookTabs(0, 2) | Out-Null
